{"js": "// The document stores the text \"<id>p028r_1</id>\" split across three\n// separate runs (tag-colored \"<id>\", plain \"p028r_1\", tag-colored\n// \"</id>\"). The edit collapses those three runs into a single run\n// containing the full \"<id>p028r_1</id>\" text, using the formatting\n// of the first (\"<id>\") run for the whole thing.\n//\n// Word's search() operates on the paragraph's logical text (it is not\n// fooled by run boundaries), so searching for the full tag text reliably\n// finds this exact paragraph (the sibling \"<id>fig_p028r_1</id>\" paragraph\n// does NOT match, since the inner text differs).\nconst results = context.document.body.search(\"<id>p028r_1</id>\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target text '<id>p028r_1</id>' in the document.\");\n}\n\n// Replacing the found range's text merges the underlying runs into one,\n// and the merged run takes on the formatting of the first of the\n// replaced runs (matching the target \"<id>\" run's Courier New / 7f6000\n// / 18 half-point formatting).\nconst target = results.items[0];\ntarget.insertText(\"<id>p028r_1</id>\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The document stores the text \"<id>p028r_1</id>\" split across three\n# separate runs (tag-colored \"<id>\", plain \"p028r_1\", tag-colored\n# \"</id>\"). The edit collapses those three runs into a single run\n# containing the full \"<id>p028r_1</id>\" text, using the formatting of\n# the first (\"<id>\") run for the whole thing.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"<id>p028r_1</id>\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"<id>p028r_1</id>\"\n$find.Forward = $true\n$find.Wrap = 0            # wdFindStop - don't wrap back over the document\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceOne (1) below is passed as the Replace parameter so only the\n# single match is touched (the sibling \"<id>fig_p028r_1</id>\" paragraph\n# has different inner text and will not match this search string).\n$found = $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, `\n    $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, `\n    $find.Replacement.Text, 1)\n\nif (-not $found) {\n    throw \"Could not find target text '<id>p028r_1</id>' in the document.\"\n}\n"}
